$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each row 2..51, column D holds a monthly mean that must be divided by the
# number of weeks belonging to that month (column B), rounded to 3 decimals.
# Column F (weekly_share) is recomputed as E / (original D / week-count), rounded
# to 3 decimals.

$lastRow = 51

for ($r = 2; $r -le $lastRow; $r++) {
    $month = $ws.Cells.Item($r, 2).Value2
    if ($month -eq $null) { continue }

    # Count how many rows (weeks) share this month value.
    $count = 0
    for ($rr = 2; $rr -le $lastRow; $rr++) {
        if ($ws.Cells.Item($rr, 2).Value2 -eq $month) {
            $count = $count + 1
        }
    }

    $oldD = $ws.Cells.Item($r, 4).Value2
    $E = $ws.Cells.Item($r, 5).Value2

    $unroundedD = $oldD / $count
    $newD = $excel.WorksheetFunction.Round($unroundedD, 3)
    $newF = $excel.WorksheetFunction.Round($E / $unroundedD, 3)

    $ws.Cells.Item($r, 4).Value2 = $newD
    $ws.Cells.Item($r, 6).Value2 = $newF
}
